# Oklahoma_B team-specific transition matrix refresh: "added more games,
# sped up simulate game logic, and drafted optimization logic" changed the
# underlying game counts feeding these transition probabilities, so the
# observed-frequency values across several states (rows) are updated below.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 0.1246819338422392
    "C2" = 0.6717557251908397
    "J2" = 0.01526717557251908
    "P2" = 0.1043256997455471
    "S2" = 0.08396946564885496
    "C3" = 0.03985507246376811
    "P3" = 0.7862318840579711
    "S3" = 0.1739130434782609
    "J4" = 0.05555555555555555
    "O4" = 0.01388888888888889
    "P4" = 0.6944444444444444
    "S4" = 0.2361111111111111
    "B6" = 0.06182795698924731
    "D6" = 0.01344086021505376
    "E6" = 0.002688172043010753
    "F6" = 0.0456989247311828
    "J6" = 0.2365591397849462
    "O6" = 0.02956989247311828
    "Q6" = 0.1666666666666667
    "R6" = 0.05376344086021505
    "S6" = 0.3897849462365591
    "B7" = 0.1111111111111111
    "D7" = 0.02116402116402116
    "E7" = 0.002645502645502645
    "F7" = 0.06084656084656084
    "J7" = 0.1005291005291005
    "O7" = 0.02645502645502645
    "Q7" = 0.2037037037037037
    "R7" = 0.06349206349206349
    "S7" = 0.4100529100529101
    "B8" = 0.09664948453608248
    "D8" = 0.01804123711340206
    "E8" = 0.001288659793814433
    "F8" = 0.07860824742268041
    "J8" = 0.05541237113402062
    "O8" = 0.02319587628865979
    "Q8" = 0.1842783505154639
    "R8" = 0.06572164948453608
    "S8" = 0.4768041237113402
    "B9" = 0.0707070707070707
    "D9" = 0.0202020202020202
    "F9" = 0.04040404040404041
    "J9" = 0.06818181818181818
    "O9" = 0.0202020202020202
    "Q9" = 0.1843434343434343
    "R9" = 0.1035353535353535
    "S9" = 0.4924242424242424
    "B10" = 0.09710494571773221
    "D10" = 0.01990349819059107
    "E10" = 0.0006031363088057901
    "F10" = 0.05729794933655006
    "J10" = 0.09589867310012062
    "O10" = 0.01930036188178528
    "Q10" = 0.2273823884197829
    "R10" = 0.07418576598311219
    "S10" = 0.4083232810615199
    "F11" = 0.001692047377326565
    "G11" = 0.1404399323181049
    "J11" = 0.06429780033840947
    "K11" = 0.2131979695431472
    "L11" = 0.5126903553299492
    "S11" = 0.06768189509306261
    "G12" = 0.7461300309597523
    "J12" = 0.07430340557275542
    "K12" = 0.01238390092879257
    "L12" = 0.04024767801857585
    "S12" = 0.1269349845201238
    "G13" = 0.6973684210526315
    "J13" = 0.1447368421052632
    "S13" = 0.1578947368421053
    "F15" = 0.04207920792079208
    "H15" = 0.1163366336633663
    "I15" = 0.08415841584158416
    "J15" = 0.25
    "K15" = 0.06930693069306931
    "M15" = 0.009900990099009901
    "N15" = 0.004950495049504951
    "O15" = 0.07425742574257425
    "S15" = 0.349009900990099
    "F16" = 0.02727272727272727
    "H16" = 0.1606060606060606
    "I16" = 0.08484848484848485
    "J16" = 0.3121212121212121
    "K16" = 0.08787878787878788
    "M16" = 0.02727272727272727
    "N16" = 0.00303030303030303
    "O16" = 0.06666666666666667
    "S16" = 0.2303030303030303
    "F17" = 0.02
    "H17" = 0.1682352941176471
    "I17" = 0.1035294117647059
    "J17" = 0.3364705882352941
    "K17" = 0.1070588235294118
    "M17" = 0.01529411764705882
    "N17" = 0.001176470588235294
    "O17" = 0.05764705882352941
    "S17" = 0.1905882352941176
    "F18" = 0.01360544217687075
    "H18" = 0.1598639455782313
    "I18" = 0.1258503401360544
    "J18" = 0.3299319727891156
    "K18" = 0.119047619047619
    "M18" = 0.02040816326530612
    "O18" = 0.04761904761904762
    "S18" = 0.1836734693877551
    "F19" = 0.02495201535508637
    "H19" = 0.1915547024952015
    "I19" = 0.08061420345489444
    "J19" = 0.2533589251439539
    "K19" = 0.09827255278310941
    "N19" = 0.0007677543186180423
    "O19" = 0.06641074856046066
    "S19" = 0.2652591170825336
    "M19" = 0.018809980806142036
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}